# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" header on the "Weekly Quantity" sheet
#    to "Weekly_PO_Qty".
# 2) Rename the "Requested quantity" header on the "Monthly Trend" sheet
#    to "Monthly_PO_Qty".
# 3) Add a new "PO Forecast" sheet (after "Monthly Trend") containing the
#    PO forecast (ds / PO_Forecast / yhat_lower / yhat_upper) data.

$wb = $excel.ActiveWorkbook

# --- 1) Weekly Quantity sheet -------------------------------------------------
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Monthly Trend sheet ---------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) New "PO Forecast" sheet -----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Header row - copy formatting (bold/border/centered) from an existing header.
$weekly.Range("A1:B1").Copy($forecast.Range("A1:B1"))
$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$weekly.Range("A1:B1").Copy($forecast.Range("C1:D1"))
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# Data rows - copy the date-formatted style from the weekly sheet's date column
# for column A, then fill in the forecast values.
$weekly.Range("A2").Copy($forecast.Range("A2:A15"))

$data = @(
    @(45060.99999999999, 10, -12.09681362973539,   30.92011654607604),
    @(45081.99999999999, 11, -9.250588065546905,   31.36302263337923),
    @(45130.99999999999, 15, -6.129745660239054,   32.9538926303206),
    @(45214.99999999999, 21, -0.7003390815725787,  40.23817693847901),
    @(45221.99999999999, 21, 2.103286479288307,    42.85050981623527),
    @(45242.99999999999, 23, 2.742914231982568,    43.88509289316186),
    @(45249.99999999999, 23, 2.535425729639869,    43.1682354332771),
    @(45256.99999999999, 24, 2.888682834820028,    43.20037508186047),
    @(45263.99999999999, 24, 4.337094875083586,    43.56419572980332),
    @(45270.99999999999, 25, 4.179925205010953,    44.43936499521243),
    @(45277.99999999999, 25, 4.545433102863711,    42.97221773243222),
    @(45284.99999999999, 25, 5.064194334511583,    44.41313639789908),
    @(45291.99999999999, 26, 6.145694317622994,    46.65944295591617),
    @(45298.99999999999, 26, 7.825972936921734,    47.64627554835117)
)

$r = 2
foreach ($row in $data) {
    $forecast.Cells.Item($r, 1).Value = $row[0]
    $forecast.Cells.Item($r, 2).Value = $row[1]
    $forecast.Cells.Item($r, 3).Value = $row[2]
    $forecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
